$wb = $excel.ActiveWorkbook

# Rename the first sheet "yli_Proteomics_Detected" -> "yli_Prots_Detected"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "yli_Prots_Detected"

# Make this sheet the active tab, and set its selection to E26
$ws1.Activate()
$ws1.Range("E26").Select()
